# Bump the "Förändrad" (Changed) date in column C by one day (45637 -> 45638)
# for every data row (rows 2 through 34) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45637) {
        $cell.Value2 = 45638
    }
}
